$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.671.13'
$ws.Range("E2").Value = '  +5.59%  '
$ws.Range("D3").Value = '2.735.15'
$ws.Range("E3").Value = '  +4.66%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''592.69'
$ws.Range("E5").Value = '  +1.42%  '
$ws.Range("D6").Value = '''153.18'
$ws.Range("E6").Value = '  +6.92%  '
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("E8").Value = '  +2.17%  '
$ws.Range("D9").Value = '2.768.24'
$ws.Range("E9").Value = '  +5.54%  '
$ws.Range("E10").Value = '  +3.55%  '
$ws.Range("E11").Value = '  +7.36%  '
$ws.Range("E12").Value = '  +3.57%  '
$ws.Range("D13").Value = '''0.159'
$ws.Range("E13").Value = '  +1.97%  '
$ws.Range("D14").Value = '3.224.32'
$ws.Range("E14").Value = '  +4.91%  '
$ws.Range("D15").Value = '''26.61'
$ws.Range("E15").Value = '  +6.63%  '
$ws.Range("D16").Value = '63.599.88'
$ws.Range("E16").Value = '  +5.46%  '
$ws.Range("E17").Value = '  +8.68%  '
$ws.Range("D18").Value = '2.758.95'
$ws.Range("E18").Value = '  +5.44%  '
$ws.Range("D19").Value = '''12.07'
$ws.Range("E20").Value = '  +4.57%  '
$ws.Range("D21").Value = '''365.62'
$ws.Range("E21").Value = '  +5.30%  '
$ws.Range("E22").Value = '  +1.59%  '
$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").Value = '''0.537'
$ws.Range("E23").Value = '  +0.80%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '''0.994'
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("D25").Value = '''65.98'
$ws.Range("E25").Value = '  +3.55%  '
$ws.Range("E26").Value = '  +4.95%  '
$ws.Range("D27").Value = '''8.66'
$ws.Range("E27").Value = '  +7.93%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").Value = '0.0₃0909'
$ws.Range("E29").Value = '  +13.59%  '
$ws.Range("E30").Value = '  +4.67%  '
$ws.Range("D31").Value = '''7.10'
$ws.Range("E31").Value = '  +9.35%  '
$ws.Range("D32").Value = '''172.86'
$ws.Range("E32").Value = '  +2.44%  '
$ws.Range("E33").Value = '  +18.43%  '
$ws.Range("E34").Value = '  -0.18%  '
$ws.Range("E35").Value = '  +5.83%  '
$ws.Range("E36").Value = '  +12.74%  '
$ws.Range("E37").Value = '  +9.29%  '
$ws.Range("E38").Value = '  +9.03%  '
$ws.Range("E39").Value = '  +18.97%  '
$ws.Range("D40").Value = '''348.40'
$ws.Range("E40").Value = '  +9.03%  '
$ws.Range("D41").Value = '''4.23'
$ws.Range("E41").Value = '  +7.34%  '
$ws.Range("D42").Value = '''38.94'
$ws.Range("E42").Value = '  +1.28%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''5.63'
$ws.Range("E43").Value = '  +11.94%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '''22.12'
$ws.Range("E44").Value = '  +10.55%  '
$ws.Range("D45").Value = '''143.51'
$ws.Range("E45").Value = '  +5.92%  '
$ws.Range("D46").Value = '''22.22'
$ws.Range("E46").Value = '  +10.84%  '
$ws.Range("E47").Value = '  +7.11%  '
$ws.Range("D48").Value = '''0.648'
$ws.Range("E48").Value = '  +6.38%  '
$ws.Range("E49").Value = '  +7.53%  '
$ws.Range("E50").Value = '  +2.71%  '
$ws.Range("D51").Value = '2.177.60'
$ws.Range("E51").Value = '  +7.51%  '
